$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "'004"
$ws.Range("N2").Value = "2019-09-30 00:00:00"
$ws.Range("O2").Value = 111398480.1
$ws.Range("P2").Value = 1257750471.36
$ws.Range("Q2").Value = 1134168326.99
$ws.Range("R2").Value = ""
$ws.Range("S2").Value = 1033846048.8
$ws.Range("T2").Value = 1033846048.8
$ws.Range("U2").Value = ""
$ws.Range("V2").Value = 28350906.01
$ws.Range("W2").Value = 13263087.91
$ws.Range("X2").Value = 17156298.84
$ws.Range("Y2").Value = 132563315.45
$ws.Range("Z2").Value = 132429619.09
$ws.Range("AA2").Value = 21030275.44
$ws.Range("AG2").Value = 5378883.03
$ws.Range("AP2").Value = ""
$ws.Range("AQ2").Value = ""
$ws.Range("AR2").Value = ""
$ws.Range("AS2").Value = 108158280.1
$ws.Range("AT2").Value = ""
